# Update "想去人数" (want-to-go count, column F) figures on the two sheets
# that carry the full event listing: "展览" (sheet 1) and "全部类型" (sheet 4).
# The row offsets differ by one between the two sheets, but each pair of
# rows refers to the same event.

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> row : new F value
$sheetA = $wb.Worksheets.Item("展览")
$sheetA.Range("F5").Value  = 5025
$sheetA.Range("F8").Value  = 13
$sheetA.Range("F9").Value  = 555
$sheetA.Range("F10").Value = 515
$sheetA.Range("F11").Value = 1038
$sheetA.Range("F13").Value = 1402
$sheetA.Range("F14").Value = 3687
$sheetA.Range("F16").Value = 137
$sheetA.Range("F17").Value = 124
$sheetA.Range("F19").Value = 2694
$sheetA.Range("F20").Value = 132
$sheetA.Range("F21").Value = 19
$sheetA.Range("F22").Value = 89
$sheetA.Range("F23").Value = 37
$sheetA.Range("F25").Value = 59
$sheetA.Range("F28").Value = 272
$sheetA.Range("F30").Value = 2

# Sheet "全部类型" -> row : new F value
$sheetB = $wb.Worksheets.Item("全部类型")
$sheetB.Range("F6").Value  = 5025
$sheetB.Range("F9").Value  = 13
$sheetB.Range("F10").Value = 555
$sheetB.Range("F11").Value = 515
$sheetB.Range("F12").Value = 1038
$sheetB.Range("F14").Value = 1402
$sheetB.Range("F15").Value = 3687
$sheetB.Range("F17").Value = 137
$sheetB.Range("F18").Value = 124
$sheetB.Range("F20").Value = 2694
$sheetB.Range("F21").Value = 132
$sheetB.Range("F22").Value = 19
$sheetB.Range("F23").Value = 89
$sheetB.Range("F24").Value = 37
$sheetB.Range("F26").Value = 59
$sheetB.Range("F29").Value = 272
$sheetB.Range("F31").Value = 2
